$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (safe from numeric auto-conversion)
$ws.Range("D2").Value = "26.555.01"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.728.69"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "1.731.43"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "26.547.35"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "1.951.74"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -3.68%  "
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("E38").Value = "  +4.29%  "
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E41").Value = "  -10.58%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  -1.31%  "

# Numeric-looking text values that must be preserved exactly as text
# (e.g. trailing zeros, multi-dot thousand separators) -
# temporarily force text format, set value, then restore default style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4802"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2665"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06168"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07179"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6077"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.521"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006952"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.517"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.787"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.243"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.776"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.405"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.970"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08004"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.686"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04504"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9073"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.053"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.405"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01501"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.499"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3883"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05378"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.850"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.246"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3405"
$ws.Range("D51").Style = "Normal"
